# scrum.xlsx — "issues 9 and 10 still pending"
#
# Append a new pending issue (#10) to the tracker: a new row recording that
# chords should be an octave lower than the melody, flagged in red like the
# other still-pending issues, then leave the selection where the user's
# cursor ended up after entering the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 10: date / issue text / status, matching the existing table layout
# (col A = date, col B = issue, col C = status, col D = notes [unused here]).
$ws.Range("A10").Value = 42599
$ws.Range("A10").NumberFormat = $ws.Range("A9").NumberFormat

$ws.Range("B10").Value = "chords should be an octave lowerr than the melody "
$ws.Range("B10").Interior.Color = 255

$ws.Range("C10").Value = "Pending"

# Leave the selection below the newly entered data, same as the author did.
$ws.Range("B12").Select()
